# Update the "想去人数" (F column) figures on both the "展览" sheet and the
# "全部类型" sheet (which mirrors the exhibition rows) to reflect the newly
# generated output.

$wb = $excel.ActiveWorkbook

$newValues = @{
    2 = 596
    3 = 59
    4 = 30
    5 = 18
    6 = 43
    7 = 35
    8 = 564
    9 = 3740
    10 = 68
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Cells.Item($row, 6).Value = $newValues[$row]
    }
}
